$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "233.44", "1.001",
# "0.4650") but must stay as literal text (trailing zeros, multi-dot
# "thousands" groupings like "30.108.92" must be preserved verbatim, and
# leading zeros/precision must not be renormalized by numeric coercion).
# Setting NumberFormat to "@" (Text) before assigning Value forces Excel to
# store the value verbatim as a string instead of inferring a number; the
# Style is then reset back to "Normal" so the cell doesn't end up carrying
# a stray Text-format style that wasn't part of the intended edit.
function Set-TextValue($cell, $value) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $value
    $ws.Range($cell).Style = "Normal"
}

$rowUpdates = @{
    2  = @{ D = "30.108.92";    E = "  -2.27%  " }
    3  = @{ D = "1.856.30";     E = "  -3.65%  " }
    4  = @{               E = "  +0.20%  " }
    5  = @{ D = "233.44";       E = "  -3.34%  " }
    6  = @{               E = "  +0.16%  " }
    7  = @{ D = "0.4650";       E = "  -2.76%  " }
    8  = @{               E = "  -2.34%  " }
    9  = @{ D = "0.06546";      E = "  -3.57%  " }
    10 = @{ D = "19.74";        E = "  +0.48%  " }
    11 = @{ D = "0.07811";      E = "  +0.22%  " }
    12 = @{ D = "96.68";        E = "  -7.21%  " }
    13 = @{ D = "1.864.83";     E = "  -3.23%  " }
    14 = @{ D = "5.102";        E = "  -3.38%  " }
    15 = @{ D = "0.6638";       E = "  -2.74%  " }
    16 = @{ D = "281.15";       E = "  -3.77%  " }
    17 = @{ D = "30.140.18";    E = "  -2.15%  " }
    18 = @{ D = "1.001";        E = "  +0.12%  " }
    19 = @{ D = "5.436";        E = "  -1.37%  " }
    20 = @{ D = "12.57";        E = "  -2.37%  " }
    21 = @{ D = "2.110.99";     E = "  -3.07%  " }
    22 = @{ D = "0.000007221"; E = "  -4.80%  " }
    23 = @{ D = "1.002";        E = "  +0.25%  " }
    24 = @{ D = "6.123";        E = "  -4.11%  " }
    25 = @{ D = "167.78";       E = "  -0.17%  " }
    26 = @{ D = "9.313";        E = "  -2.36%  " }
    27 = @{ D = "18.84";        E = "  -4.75%  " }
    28 = @{ D = "1.909";        E = "  -9.77%  " }
    29 = @{ D = "1.335";        E = "  -4.06%  " }
    30 = @{ D = "0.09545";      E = "  -5.35%  " }
    31 = @{ D = "4.413";        E = "  -4.19%  " }
    32 = @{ D = "1.468";        E = "  -3.88%  " }
    33 = @{ D = "4.095";        E = "  -5.35%  " }
    34 = @{ D = "0.04640";      E = "  -3.65%  " }
    35 = @{ B = "ARBITRUM";   C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";     D = "1.097";  E = "  -2.50%  " }
    36 = @{ B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "0.6990"; E = "  -4.92%  " }
    37 = @{ D = "2.701";        E = "  -0.58%  " }
    38 = @{ D = "0.01843";      E = "  -5.13%  " }
    39 = @{               E = "  -1.45%  " }
    40 = @{ D = "2.512";        E = "  -4.55%  " }
    41 = @{ D = "71.96";        E = "  -4.11%  " }
    42 = @{ D = "0.8520";       E = "  -1.83%  " }
    43 = @{ D = "1.914";        E = "  -5.34%  " }
    44 = @{ D = "1.001";        E = "  +0.09%  " }
    45 = @{ D = "103.94";       E = "  -1.80%  " }
    46 = @{ D = "0.4136";       E = "  -4.75%  " }
    47 = @{ D = "988.80";       E = "  +0.09%  " }
    48 = @{ D = "7.189";        E = "  -4.61%  " }
    49 = @{ D = "9.149";        E = "  +1.59%  " }
    50 = @{ D = "33.92";        E = "  -2.90%  " }
    51 = @{ D = "0.1137";       E = "  -6.19%  " }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cols = $rowUpdates[$rowNum]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$rowNum"
        $value = $cols[$col]
        if ($col -eq "D") {
            Set-TextValue $cellRef $value
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
